$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Insert a new row above row 2 (shifts existing data + totals rows down by one)
$ws.Rows.Item(2).Insert()

# Fill in the new row 2 with the latest day's data (BRS daily figures)
$ws.Cells.Item(2, 1).Value = 45446
$ws.Cells.Item(2, 2).Value = 212
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 211

# Copy the date number format from the row below onto the new date cell
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update the totals row formulas to cover the newly widened data range (rows 2:4)
$ws.Cells.Item(5, 2).Formula = "=SUM(B2:B4)"
$ws.Range("C5:E5").Formula = "=SUM(C2:C4)"

$ws.Range("E7").Select()
